$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Robo1"
$ws.Cells.Item(2, 3).Value = "Robo1"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1325553333333333
$ws.Cells.Item(2, 8).Value = 0.397666
$ws.Cells.Item(2, 9).Value = 0.00533964316398423
$ws.Cells.Item(2, 10).Value = 0.00533964316398423
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.1325553333333333
$ws.Cells.Item(2, 14).Value = 0.397666
$ws.Cells.Item(2, 15).Value = 0.00533964316398423
$ws.Cells.Item(2, 16).Value = 0.00533964316398423
$ws.Cells.Item(2, 17).Value = 0.01757091639511111
$ws.Cells.Item(2, 18).Value = 0.158138247556
$ws.Cells.Item(2, 19).Value = 0.00002851178911868352
$ws.Cells.Item(2, 20).Value = 0.00002851178911868352

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Robo1"
$ws.Cells.Item(3, 3).Value = "Robo1"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1325553333333333
$ws.Cells.Item(3, 8).Value = 0.397666
$ws.Cells.Item(3, 9).Value = 0.00533964316398423
$ws.Cells.Item(3, 10).Value = 0.00533964316398423
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 17.178266
$ws.Cells.Item(3, 14).Value = 51.534798
$ws.Cells.Item(3, 15).Value = 0.6919812904497951
$ws.Cells.Item(3, 16).Value = 0.691981290449795
$ws.Cells.Item(3, 17).Value = 2.277070775718667
$ws.Cells.Item(3, 18).Value = 20.493636981468
$ws.Cells.Item(3, 19).Value = 0.003694933167155234
$ws.Cells.Item(3, 20).Value = 0.003694933167155233

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Robo1"
$ws.Cells.Item(4, 3).Value = "Robo1"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1325553333333333
$ws.Cells.Item(4, 8).Value = 0.397666
$ws.Cells.Item(4, 9).Value = 0.00533964316398423
$ws.Cells.Item(4, 10).Value = 0.00533964316398423
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 7.513933666666667
$ws.Cells.Item(4, 14).Value = 22.541801
$ws.Cells.Item(4, 15).Value = 0.3026790663862208
$ws.Cells.Item(4, 16).Value = 0.3026790663862208
$ws.Cells.Item(4, 17).Value = 0.9960119818295555
$ws.Cells.Item(4, 18).Value = 8.964107836466
$ws.Cells.Item(4, 19).Value = 0.001616198207710313
$ws.Cells.Item(4, 20).Value = 0.001616198207710313

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Robo1"
$ws.Cells.Item(5, 3).Value = "Robo1"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 17.178266
$ws.Cells.Item(5, 8).Value = 51.534798
$ws.Cells.Item(5, 9).Value = 0.6919812904497951
$ws.Cells.Item(5, 10).Value = 0.691981290449795
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.1325553333333333
$ws.Cells.Item(5, 14).Value = 0.397666
$ws.Cells.Item(5, 15).Value = 0.00533964316398423
$ws.Cells.Item(5, 16).Value = 0.00533964316398423
$ws.Cells.Item(5, 17).Value = 2.277070775718667
$ws.Cells.Item(5, 18).Value = 20.493636981468
$ws.Cells.Item(5, 19).Value = 0.003694933167155234
$ws.Cells.Item(5, 20).Value = 0.003694933167155233

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Robo1"
$ws.Cells.Item(6, 3).Value = "Robo1"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 17.178266
$ws.Cells.Item(6, 8).Value = 51.534798
$ws.Cells.Item(6, 9).Value = 0.6919812904497951
$ws.Cells.Item(6, 10).Value = 0.691981290449795
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 17.178266
$ws.Cells.Item(6, 14).Value = 51.534798
$ws.Cells.Item(6, 15).Value = 0.6919812904497951
$ws.Cells.Item(6, 16).Value = 0.691981290449795
$ws.Cells.Item(6, 17).Value = 295.092822766756
$ws.Cells.Item(6, 18).Value = 2655.835404900804
$ws.Cells.Item(6, 19).Value = 0.4788381063325636
$ws.Cells.Item(6, 20).Value = 0.4788381063325635

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Robo1"
$ws.Cells.Item(7, 3).Value = "Robo1"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 17.178266
$ws.Cells.Item(7, 8).Value = 51.534798
$ws.Cells.Item(7, 9).Value = 0.6919812904497951
$ws.Cells.Item(7, 10).Value = 0.691981290449795
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 7.513933666666667
$ws.Cells.Item(7, 14).Value = 22.541801
$ws.Cells.Item(7, 15).Value = 0.3026790663862208
$ws.Cells.Item(7, 16).Value = 0.3026790663862208
$ws.Cells.Item(7, 17).Value = 129.0763512323553
$ws.Cells.Item(7, 18).Value = 1161.687161091198
$ws.Cells.Item(7, 19).Value = 0.2094482509500762
$ws.Cells.Item(7, 20).Value = 0.2094482509500762

# Row 8
$ws.Cells.Item(8, 1).Value = "ECs"
$ws.Cells.Item(8, 2).Value = "Robo1"
$ws.Cells.Item(8, 3).Value = "Robo1"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 7.513933666666667
$ws.Cells.Item(8, 8).Value = 22.541801
$ws.Cells.Item(8, 9).Value = 0.3026790663862208
$ws.Cells.Item(8, 10).Value = 0.3026790663862208
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.1325553333333333
$ws.Cells.Item(8, 14).Value = 0.397666
$ws.Cells.Item(8, 15).Value = 0.00533964316398423
$ws.Cells.Item(8, 16).Value = 0.00533964316398423
$ws.Cells.Item(8, 17).Value = 0.9960119818295555
$ws.Cells.Item(8, 18).Value = 8.964107836466
$ws.Cells.Item(8, 19).Value = 0.001616198207710313
$ws.Cells.Item(8, 20).Value = 0.001616198207710313

# Row 9
$ws.Cells.Item(9, 1).Value = "ECs"
$ws.Cells.Item(9, 2).Value = "Robo1"
$ws.Cells.Item(9, 3).Value = "Robo1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 7.513933666666667
$ws.Cells.Item(9, 8).Value = 22.541801
$ws.Cells.Item(9, 9).Value = 0.3026790663862208
$ws.Cells.Item(9, 10).Value = 0.3026790663862208
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 17.178266
$ws.Cells.Item(9, 14).Value = 51.534798
$ws.Cells.Item(9, 15).Value = 0.6919812904497951
$ws.Cells.Item(9, 16).Value = 0.691981290449795
$ws.Cells.Item(9, 17).Value = 129.0763512323553
$ws.Cells.Item(9, 18).Value = 1161.687161091198
$ws.Cells.Item(9, 19).Value = 0.2094482509500762
$ws.Cells.Item(9, 20).Value = 0.2094482509500762

# Row 10
$ws.Cells.Item(10, 1).Value = "ECs"
$ws.Cells.Item(10, 2).Value = "Robo1"
$ws.Cells.Item(10, 3).Value = "Robo1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 7.513933666666667
$ws.Cells.Item(10, 8).Value = 22.541801
$ws.Cells.Item(10, 9).Value = 0.3026790663862208
$ws.Cells.Item(10, 10).Value = 0.3026790663862208
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 7.513933666666667
$ws.Cells.Item(10, 14).Value = 22.541801
$ws.Cells.Item(10, 15).Value = 0.3026790663862208
$ws.Cells.Item(10, 16).Value = 0.3026790663862208
$ws.Cells.Item(10, 17).Value = 56.45919914706678
$ws.Cells.Item(10, 18).Value = 508.132792323601
$ws.Cells.Item(10, 19).Value = 0.09161461722843424
$ws.Cells.Item(10, 20).Value = 0.09161461722843424
